$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 11:01:43"
$ws.Range("A3").Value = "Total filas: 184"
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "07:18:07"
$arr[0,1] = "07:21"
$arr[0,2] = "10_OLMOS"
$arr[0,3] = 3
$arr[0,4] = "LP1912"
$ws.Range("A47:E47").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "07:18:07"
$arr[0,1] = "07:21"
$arr[0,2] = "225_GOMEZ"
$arr[0,3] = 3
$arr[0,4] = "LP1912"
$ws.Range("A48:E48").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "07:18:07"
$arr[0,1] = "07:21"
$arr[0,2] = "26_HERNANDEZ"
$arr[0,3] = 3
$arr[0,4] = "LP1912"
$ws.Range("A49:E49").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "08:16:28"
$arr[0,1] = "08:23"
$arr[0,2] = "215B_EL PATO"
$arr[0,3] = 7
$arr[0,4] = "LP1912"
$ws.Range("A78:E78").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "07:18:07"
$arr[0,1] = "08:23"
$arr[0,2] = "16_P MOR-SANTA ANA"
$arr[0,3] = 65
$arr[0,4] = "LP1912"
$ws.Range("A79:E79").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "08:38:27"
$arr[0,1] = "09:23"
$arr[0,2] = "17_ROMERO"
$arr[0,3] = 45
$arr[0,4] = "LP1912"
$ws.Range("A108:E108").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "08:52:20"
$arr[0,1] = "09:23"
$arr[0,2] = "16_SANTA ANA"
$arr[0,3] = 31
$arr[0,4] = "LP1912"
$ws.Range("A110:E110").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:01"
$arr[0,2] = "11_ETCHEVERRY"
$arr[0,3] = 0
$arr[0,4] = "LP1912"
$ws.Range("A148:E148").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:01"
$arr[0,2] = "215C_EL PATO"
$arr[0,3] = 0
$arr[0,4] = "LP1912"
$ws.Range("A149:E149").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:01"
$arr[0,2] = "16_SANTA ANA"
$arr[0,3] = 0
$arr[0,4] = "LP1912"
$ws.Range("A150:E150").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "11:02"
$arr[0,2] = "215C_EL PATO"
$arr[0,3] = 35
$arr[0,4] = "LP1912"
$ws.Range("A151:E151").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "11:04"
$arr[0,2] = "11_ETCHEVERRY"
$arr[0,3] = 37
$arr[0,4] = "LP1912"
$ws.Range("A152:E152").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:05"
$arr[0,2] = "23_HERNANDEZ"
$arr[0,3] = 4
$arr[0,4] = "LP1912"
$ws.Range("A153:E153").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:06"
$arr[0,2] = "16_P MOR-167 Y 521"
$arr[0,3] = 5
$arr[0,4] = "LP1912"
$ws.Range("A154:E154").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "11:07"
$arr[0,2] = "16_P MOR-167 Y 521"
$arr[0,3] = 40
$arr[0,4] = "LP1912"
$ws.Range("A155:E155").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:11"
$arr[0,2] = "10_OLMOS"
$arr[0,3] = 10
$arr[0,4] = "LP1912"
$ws.Range("A156:E156").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:12"
$arr[0,2] = "15_ABASTO"
$arr[0,3] = 11
$arr[0,4] = "LP1912"
$ws.Range("A157:E157").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "11:17"
$arr[0,2] = "16_SANTA ANA"
$arr[0,3] = 50
$arr[0,4] = "LP1912"
$ws.Range("A158:E158").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:19"
$arr[0,2] = "86_EST CHICA-ESC AGRARIA"
$arr[0,3] = 18
$arr[0,4] = "LP1912"
$ws.Range("A159:E159").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "11:20"
$arr[0,2] = "86_EST CHICA-ESC AGRARIA"
$arr[0,3] = 53
$arr[0,4] = "LP1912"
$ws.Range("A160:E160").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "09:27:56"
$arr[0,1] = "11:20"
$arr[0,2] = "26_HERNANDEZ"
$arr[0,3] = 113
$arr[0,4] = "LP1912"
$ws.Range("A161:E161").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:21"
$arr[0,2] = "26_HERNANDEZ"
$arr[0,3] = 20
$arr[0,4] = "LP1912"
$ws.Range("A162:E162").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:26"
$arr[0,2] = "16_SANTA ANA"
$arr[0,3] = 25
$arr[0,4] = "LP1912"
$ws.Range("A163:E163").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:27"
$arr[0,2] = "225_C ROCA-H SUR"
$arr[0,3] = 26
$arr[0,4] = "LP1912"
$ws.Range("A164:E164").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:32"
$arr[0,2] = "81_EL PELIGRO"
$arr[0,3] = 31
$arr[0,4] = "LP1912"
$ws.Range("A165:E165").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:34"
$arr[0,2] = "23_HERNANDEZ"
$arr[0,3] = 33
$arr[0,4] = "LP1912"
$ws.Range("A166:E166").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:35"
$arr[0,2] = "11_ETCHEVERRY"
$arr[0,3] = 34
$arr[0,4] = "LP1912"
$ws.Range("A167:E167").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "11:36"
$arr[0,2] = "11_ETCHEVERRY"
$arr[0,3] = 69
$arr[0,4] = "LP1912"
$ws.Range("A168:E168").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:41"
$arr[0,2] = "17_ROMERO"
$arr[0,3] = 40
$arr[0,4] = "LP1912"
$ws.Range("A169:E169").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "11:42"
$arr[0,2] = "17_ROMERO"
$arr[0,3] = 75
$arr[0,4] = "LP1912"
$ws.Range("A170:E170").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:43"
$arr[0,2] = "10_OLMOS"
$arr[0,3] = 42
$arr[0,4] = "LP1912"
$ws.Range("A171:E171").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:45"
$arr[0,2] = "10_OLMOS"
$arr[0,3] = 44
$arr[0,4] = "LP1912"
$ws.Range("A172:E172").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "11:47"
$arr[0,2] = "10_OLMOS"
$arr[0,3] = 80
$arr[0,4] = "LP1912"
$ws.Range("A173:E173").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:51"
$arr[0,2] = "215B_EL PATO"
$arr[0,3] = 50
$arr[0,4] = "LP1912"
$ws.Range("A174:E174").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:52"
$arr[0,2] = "15_ABASTO"
$arr[0,3] = 51
$arr[0,4] = "LP1912"
$ws.Range("A175:E175").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:59"
$arr[0,2] = "225_GOMEZ"
$arr[0,3] = 58
$arr[0,4] = "LP1912"
$ws.Range("A176:E176").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:02"
$arr[0,2] = "84_COLONIA URQUIZA-ESC 49"
$arr[0,3] = 61
$arr[0,4] = "LP1912"
$ws.Range("A177:E177").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:06"
$arr[0,2] = "16_P MOR-SANTA ANA"
$arr[0,3] = 65
$arr[0,4] = "LP1912"
$ws.Range("A178:E178").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:06"
$arr[0,2] = "14_ABASTO"
$arr[0,3] = 65
$arr[0,4] = "LP1912"
$ws.Range("A179:E179").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "12:07"
$arr[0,2] = "16_P MOR-SANTA ANA"
$arr[0,3] = 100
$arr[0,4] = "LP1912"
$ws.Range("A180:E180").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:14"
$arr[0,2] = "17_ROMERO"
$arr[0,3] = 73
$arr[0,4] = "LP1912"
$ws.Range("A181:E181").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "12:19"
$arr[0,2] = "14_ABASTO"
$arr[0,3] = 112
$arr[0,4] = "LP1912"
$ws.Range("A182:E182").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:20"
$arr[0,2] = "215A_EL PATO"
$arr[0,3] = 79
$arr[0,4] = "LP1912"
$ws.Range("A183:E183").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:20"
$arr[0,2] = "14_ABASTO"
$arr[0,3] = 79
$arr[0,4] = "LP1912"
$ws.Range("A184:E184").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:21"
$arr[0,2] = "26_HERNANDEZ"
$arr[0,3] = 80
$arr[0,4] = "LP1912"
$ws.Range("A185:E185").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "12:21"
$arr[0,2] = "215A_EL PATO"
$arr[0,3] = 114
$arr[0,4] = "LP1912"
$ws.Range("A186:E186").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:36"
$arr[0,2] = "27_EL RETIRO"
$arr[0,3] = 95
$arr[0,4] = "LP1912"
$ws.Range("A187:E187").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:38"
$arr[0,2] = "17_179 Y 38"
$arr[0,3] = 97
$arr[0,4] = "LP1912"
$ws.Range("A188:E188").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:48"
$arr[0,2] = "11_ETCHEVERRY"
$arr[0,3] = 107
$arr[0,4] = "LP1912"
$ws.Range("A189:E189").Value = $arr

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 11:01:43"
$ws.Range("A3").Value = "Total filas: 20"
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:01"
$arr[0,2] = "215C_EL PATO"
$arr[0,3] = 0
$arr[0,4] = "LP1912"
$ws.Range("A21:E21").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:51"
$arr[0,2] = "215B_EL PATO"
$arr[0,3] = 50
$arr[0,4] = "LP1912"
$ws.Range("A23:E23").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:20"
$arr[0,2] = "215A_EL PATO"
$arr[0,3] = 79
$arr[0,4] = "LP1912"
$ws.Range("A24:E24").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "10:27:08"
$arr[0,1] = "12:21"
$arr[0,2] = "215A_EL PATO"
$arr[0,3] = 114
$arr[0,4] = "LP1912"
$ws.Range("A25:E25").Value = $arr

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 11:01:43"
$ws.Range("A3").Value = "Total filas: 24"
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "11:13"
$arr[0,2] = "215C_LA PLATA"
$arr[0,3] = 12
$arr[0,4] = "L6203"
$ws.Range("A26:E26").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:04"
$arr[0,2] = "215A_LA PLATA"
$arr[0,3] = 63
$arr[0,4] = "L6173"
$ws.Range("A28:E28").Value = $arr
$arr = New-Object 'object[,]' 1,5
$arr[0,0] = "11:01:43"
$arr[0,1] = "12:53"
$arr[0,2] = "215C_LA PLATA"
$arr[0,3] = 112
$arr[0,4] = "L6203"
$ws.Range("A29:E29").Value = $arr
